$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.198.73'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -0.65%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.657.50'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -1.29%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.61%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''218.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +0.51%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.5220'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -1.67%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.57%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.2663'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +0.12%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.06327'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -1.97%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''21.17'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -0.25%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.07739'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -0.71%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = '''Polkadot'
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value = '''4.436'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -1.53%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = '''WrappedEther'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = '''1.651.02'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -1.56%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''0.5477'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -2.29%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''0.0₅8249'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -2.29%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''64.98'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -1.67%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''26.215.22'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -0.62%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = '''  +0.53%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''4.684'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -3.09%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''193.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -1.21%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''10.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -1.96%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''6.122'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -4.38%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''1.008'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +0.82%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''138.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -3.05%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''0.1238'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -2.10%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''7.276'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -2.82%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''16.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -1.18%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''1.413'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -1.53%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''0.06056'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -2.38%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''1.284'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.69%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''3.559'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.41%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''3.354'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -3.09%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''1.654'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -3.14%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''0.9840'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -2.98%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  +0.39%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''2.780'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -0.08%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''0.5941'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +3.37%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.01598'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -2.14%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''5.968'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +0.53%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''0.8663'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +0.06%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''1.049.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -0.51%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''  +0.43%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''99.95'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -0.09%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''1.795.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -1.59%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = '''BabyDogeCoin'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = '''https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = '''0.0₈109'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -0.78%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = '''Aave'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''57.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +0.47%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''1.007'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  +0.54%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''8.148'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -0.15%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''0.05187'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -0.12%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''1.477'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  +3.72%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.4231'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +0.30%  '
$ws.Range('E51').Style = 'Normal'
